# Applies the "removed reference to CHA's name and area" edit to the
# survey sheet of the MOH 515 (Post Outbreak) XLSForm.
#
# Summary of changes to the `survey` sheet:
#   1. Insert a new `hidden_inputs` group (3 rows) at the very top of the
#      question list (right after the header row).
#   2. Give the `form_summary` group (and every other top-level group) an
#      "appearence" of `field-list`.
#   3. Remove the `cha_name` ("What Is Your Name?", db:person/db-object)
#      and `cha_area` ("What Is Your Area?") questions entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

$xlPasteFormats = -4122

# --- 1. Insert 3 blank rows right after the header row (row 1) --------
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()

# The inserted rows pick up the full A:H formatting of the row that used
# to be at position 2 (which had data out to column H); trim the unused
# G:H cells on the 3 new rows back down so they stay genuinely empty.
$ws.Range("G2:H4").Clear()

# New row 2: begin group / hidden_inputs, relevant=./source='user', appearence=field-list
$ws.Cells.Item(2,1).Value = "begin group"
$ws.Cells.Item(2,2).Value = "hidden_inputs"
$ws.Cells.Item(2,3).Value = $null
$ws.Cells.Item(2,4).Value = $null
$ws.Cells.Item(2,5).Value = "./source='user'"
$ws.Cells.Item(2,6).Value = "field-list"

# New row 3: hidden / source
$ws.Cells.Item(3,1).Value = "hidden"
$ws.Cells.Item(3,2).Value = "source"

# New row 4: end group
$ws.Cells.Item(4,1).Value = "end group"

# --- 2. The (now shifted) form_summary begin-group row gains appearence
#         field-list. Old row 2 ("begin group"/form_summary) is now row 5.
#         Column F never had a cell there before, so copy the normal cell
#         format over (from D5) before writing the value, to match the
#         style ("s") used by every other populated cell.
$ws.Cells.Item(5,4).Copy()
$ws.Cells.Item(5,6).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(5,6).Value = "field-list"

# --- 3. Remove the cha_name and cha_area rows. After the insert above,
#         old row 3 (cha_name) is now row 6 and old row 4 (cha_area) is
#         now row 7 -- deleting row 6 twice removes both.
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(6).Delete()

# --- 4. The remaining top-level groups (household_indicators,
#         cholera_cases, community_activities) also gain an appearence of
#         field-list on their begin-group row (again, column F is new on
#         these rows, so copy formatting from column A first).
$ws.Cells.Item(11,1).Copy()
$ws.Cells.Item(11,6).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(11,6).Value = "field-list"

$ws.Cells.Item(18,1).Copy()
$ws.Cells.Item(18,6).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(18,6).Value = "field-list"

$ws.Cells.Item(25,1).Copy()
$ws.Cells.Item(25,6).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(25,6).Value = "field-list"
